# CRM-5138 - Removed 'Amount paid' column that was coming in downloaded file,
# also sequence of columns are same in downloaded file as on view.
#
# The template used to output 9 columns (A:I) in this order:
#   Product Name, Part Requested, Model Number, Booking Id, Age of Requested,
#   Parts Code, Parts Quantity, Amount Due, Serial Number
#
# It must become 8 columns (A:H), matching the on-screen view order and
# dropping the "Amount Due" column entirely:
#   Booking ID, Appliance, Spare Part, Parts Number, Quantity,
#   Age of Requested, Model No, Serial No

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Amount Due" column (old column H) together with its
# placeholder cell below it - this shifts the old Serial Number column (I)
# left into H, exactly like the authored edit that dropped the amount column.
$ws.Columns.Item(8).Delete()

# --- Row 1: column headers, in the new on-screen order ---
$ws.Range("A1").Value = "Booking ID"
$ws.Range("B1").Value = "Appliance"
$ws.Range("C1").Value = "Spare Part"
$ws.Range("D1").Value = "Parts Number"
$ws.Range("E1").Value = "Quantity"
$ws.Range("F1").Value = "Age of Requested"
$ws.Range("G1").Value = "Model No"
$ws.Range("H1").Value = "Serial No"

# --- Row 2: merge-field placeholders that line up with the headers above ---
$ws.Range("A2").Value = "{spare:booking_id}"
$ws.Range("B2").Value = "{spare:services}"
$ws.Range("C2").Value = "{spare:parts_requested}"
$ws.Range("D2").Value = "{spare:part_number}"
$ws.Range("E2").Value = "{spare:quantity}"
$ws.Range("F2").Value = "{spare:age_of_request}"
$ws.Range("G2").Value = "{spare:model_number}"
$ws.Range("H2").Value = "{spare:serial_number}"

# --- Column widths: keep the original look for the remaining columns ---
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(7).ColumnWidth
$ws.Columns.Item(8).ColumnWidth = 17.83
$ws.Range($ws.Cells.Item(1, 9), $ws.Cells.Item(1, 16)).EntireColumn.ColumnWidth = 7.83

# --- View: scroll back to the top-left (drop the stale topLeftCell/selection) ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("H1").Select()
